$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("B2").Value = 0.1936619718309859
    $ws.Range("C2").Value = 0.5211267605633803
    $ws.Range("J2").Value = 0.02464788732394366
    $ws.Range("P2").Value = 0.1795774647887324
    $ws.Range("S2").Value = 0.08098591549295775
    $ws.Range("B3").Value = 0.01036269430051814
    $ws.Range("C3").Value = 0.0310880829015544
    $ws.Range("J3").Value = 0.04145077720207254
    $ws.Range("P3").Value = 0.7668393782383419
    $ws.Range("S3").Value = 0.150259067357513
    $ws.Range("J4").Value = 0.1041666666666667
    $ws.Range("P4").Value = 0.625
    $ws.Range("S4").Value = 0.2708333333333333
    $ws.Range("B6").Value = 0.05517241379310345
    $ws.Range("D6").Value = 0.01379310344827586
    $ws.Range("F6").Value = 0.03448275862068965
    $ws.Range("J6").Value = 0.2620689655172414
    $ws.Range("O6").Value = 0.02068965517241379
    $ws.Range("Q6").Value = 0.1724137931034483
    $ws.Range("R6").Value = 0.08275862068965517
    $ws.Range("S6").Value = 0.3586206896551724
    $ws.Range("B7").Value = 0.1360544217687075
    $ws.Range("D7").Value = 0.006802721088435374
    $ws.Range("E7").Value = 0.006802721088435374
    $ws.Range("F7").Value = 0.0272108843537415
    $ws.Range("J7").Value = 0.1496598639455782
    $ws.Range("O7").Value = 0.01360544217687075
    $ws.Range("Q7").Value = 0.1836734693877551
    $ws.Range("R7").Value = 0.08163265306122448
    $ws.Range("S7").Value = 0.3945578231292517
    $ws.Range("B8").Value = 0.09433962264150944
    $ws.Range("D8").Value = 0.01886792452830189
    $ws.Range("F8").Value = 0.04245283018867924
    $ws.Range("J8").Value = 0.06367924528301887
    $ws.Range("O8").Value = 0.02122641509433962
    $ws.Range("Q8").Value = 0.2240566037735849
    $ws.Range("R8").Value = 0.09198113207547169
    $ws.Range("S8").Value = 0.4433962264150944
    $ws.Range("B9").Value = 0.1413612565445026
    $ws.Range("D9").Value = 0.02094240837696335
    $ws.Range("E9").Value = 0.01047120418848168
    $ws.Range("F9").Value = 0.03664921465968586
    $ws.Range("J9").Value = 0.07329842931937172
    $ws.Range("O9").Value = 0.02094240837696335
    $ws.Range("Q9").Value = 0.1989528795811518
    $ws.Range("R9").Value = 0.07329842931937172
    $ws.Range("S9").Value = 0.4240837696335079
    $ws.Range("B10").Value = 0.08708938120702826
    $ws.Range("D10").Value = 0.02444614209320092
    $ws.Range("F10").Value = 0.04889228418640183
    $ws.Range("J10").Value = 0.1558441558441558
    $ws.Range("O10").Value = 0.02750190985485103
    $ws.Range("Q10").Value = 0.2337662337662338
    $ws.Range("R10").Value = 0.07333842627960276
    $ws.Range("S10").Value = 0.3491214667685256
    $ws.Range("G11").Value = 0.1283185840707965
    $ws.Range("J11").Value = 0.1327433628318584
    $ws.Range("K11").Value = 0.1769911504424779
    $ws.Range("L11").Value = 0.5442477876106194
    $ws.Range("S11").Value = 0.01769911504424779
    $ws.Range("G12").Value = 0.7377049180327869
    $ws.Range("J12").Value = 0.2213114754098361
    $ws.Range("K12").Value = 0.00819672131147541
    $ws.Range("S12").Value = 0.03278688524590164
    $ws.Range("G13").Value = 0.725
    $ws.Range("J13").Value = 0.2
    $ws.Range("S13").Value = 0.075
    $ws.Range("F15").Value = 0.009302325581395349
    $ws.Range("H15").Value = 0.1348837209302326
    $ws.Range("I15").Value = 0.06046511627906977
    $ws.Range("J15").Value = 0.427906976744186
    $ws.Range("K15").Value = 0.05581395348837209
    $ws.Range("M15").Value = 0.01395348837209302
    $ws.Range("O15").Value = 0.06511627906976744
    $ws.Range("S15").Value = 0.2325581395348837
    $ws.Range("F16").Value = 0.009009009009009009
    $ws.Range("H16").Value = 0.1396396396396396
    $ws.Range("I16").Value = 0.06306306306306306
    $ws.Range("J16").Value = 0.490990990990991
    $ws.Range("K16").Value = 0.07207207207207207
    $ws.Range("M16").Value = 0.01801801801801802
    $ws.Range("N16").Value = 0.004504504504504504
    $ws.Range("O16").Value = 0.04954954954954955
    $ws.Range("S16").Value = 0.1531531531531531
    $ws.Range("F17").Value = 0.01727447216890595
    $ws.Range("H17").Value = 0.1593090211132438
    $ws.Range("I17").Value = 0.07101727447216891
    $ws.Range("J17").Value = 0.4203454894433781
    $ws.Range("K17").Value = 0.1151631477927063
    $ws.Range("M17").Value = 0.01919385796545105
    $ws.Range("O17").Value = 0.0671785028790787
    $ws.Range("S17").Value = 0.1305182341650672
    $ws.Range("H18").Value = 0.1005586592178771
    $ws.Range("I18").Value = 0.0893854748603352
    $ws.Range("J18").Value = 0.4636871508379888
    $ws.Range("K18").Value = 0.1229050279329609
    $ws.Range("M18").Value = 0.01675977653631285
    $ws.Range("O18").Value = 0.0893854748603352
    $ws.Range("S18").Value = 0.1173184357541899
    $ws.Range("F19").Value = 0.007563025210084034
    $ws.Range("H19").Value = 0.2285714285714286
    $ws.Range("I19").Value = 0.09411764705882353
    $ws.Range("J19").Value = 0.3815126050420168
    $ws.Range("K19").Value = 0.08907563025210084
    $ws.Range("M19").Value = 0.01680672268907563
    $ws.Range("O19").Value = 0.05798319327731093
    $ws.Range("S19").Value = 0.1275646743978591

